$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Playing Now"

# Add the new "Genre" header in column C
$ws.Range("C1").Value = "Genre"

# Set column B width
$ws.Columns.Item(2).ColumnWidth = 13.7109375

# Update selection to A2
$ws.Range("A2").Select()
